$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Score")

# --- Extend the "Cost per mile" mini table with a Price-per-gallon input -
# (done first, before the banner rows are inserted above, so the new
# "Price per gallon" shared string is allocated ahead of the banner text,
# matching the order the strings were typed in the real edit)
$ws.Range("H24").Value = $ws.Range("G24").Value2
$ws.Range("G24").Value = $ws.Range("F24").Value2
$ws.Range("F24").Value = "Price per gallon"

$ws.Range("H25").Formula = $ws.Range("G25").Formula
$ws.Range("F25").Value = 4
$ws.Range("G25").Formula = "=F25*D25/C25"
$ws.Range("G25").NumberFormat = '"$"#,##0.00'

# --- Insert the two new banner rows at the top of the sheet --------------
# Old row 3 was an empty placeholder row (A3 only, no content), so we
# insert two fresh rows above it and then drop the now-redundant old row
# (it lands one row below the two new ones after both inserts).
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(5).Delete()

$ws.Range("A3:F4").Font.Bold = $false

$ws.Range("A3").Value = "Excellent work!"
$ws.Range("A3:F3").Merge()
$ws.Range("A3:F3").WrapText = $true

$ws.Range("A4").Value = "Here's the grade breakdown:"
$ws.Range("A4:F4").Merge()
$ws.Range("A4:F4").WrapText = $true

# --- Add the wrapped comment column along the rubric rows ----------------
$ws.Range("F6:F22").WrapText = $true
